# Update countries & provincias Spain
# Applies the 6-Oct-2020 22:15 data refresh to the "Pais" sheet:
#  - refreshed case counters for a number of countries
#  - Angola overtook Suazilandia (now ranked above it)
#  - Siria overtook Bahamas (now ranked above it)
#  - Nueva Caledonia now listed above Santa Lucia (tie-break reorder)
#  - timestamp footer updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $vals)
    $cols = @("B","C","D","E","F","G","H")
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 22:15"

# Estados Unidos
Set-Row 4 @(7709109, 29465, 4920294, 2573269, 0, 514, 215546)

# India
Set-Row 5 @(6754179, 72106, 5741253, 908335, 0, 991, 104591)

# Alemania
Set-Row 26 @(307114, 2457, 263700, 33780, 0, 18, 9634)

# Israel
Set-Row 27 @(277026, 4717, 211397, 63832, 0, 40, 1797)

# Canada
Set-Row 29 @(170945, 1985, 143767, 17651, 0, 23, 9527)

# Catar
Set-Row 36 @(126943, 251, 123893, 2834)

# Costa Rica
Set-Row 51 @(82142, 1013, 50020, 31118, 0, 17, 1004)

# Argelia
Set-Row 64 @(52399, 129, 36763, 13868)

# Cabo Verde
Set-Row 118 @(6518, 85, 5632, 817, 0, 1, 69)

# Row 121 becomes Angola (overtakes Suazilandia), with its refreshed totals
$ws.Range("A121").Value = "Angola"
Set-Row 121 @(5725, 195, 2598, 2916, 0, 12, 211)

# Row 122 becomes Suazilandia, keeping its prior totals
$ws.Range("A122").Value = "Suazilandia"
Set-Row 122 @(5579, 0, 5141, 326, 0, 0, 112)

# Row 134 becomes Siria (overtakes Bahamas), with its refreshed totals
$ws.Range("A134").Value = "Siria"
Set-Row 134 @(4457, 46, 1183, 3065, 0, 2, 209)

# Row 135 becomes Bahamas, keeping its prior totals
$ws.Range("A135").Value = "Bahamas"
Set-Row 135 @(4452, 0, 2375, 1981, 0, 0, 96)

# Botsuana
Set-Row 147 @(3172, 0, 834, 2320, 0, 2, 18)

# Seychelles
Set-Row 193 @(148, 2, 143, 5, 0, 0)

# Row 207 becomes Nueva Caledonia, row 208 becomes Santa Lucia (tie-break swap, totals unchanged)
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"
